# Re-calculated grid populations for Springfield, Missouri MSA
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Notes text for the "sprfd_mo_msa_grid02" row (row 25) with the
# corrected / expanded population-interpolation description (overwritten
# with corrected calculation on 10-Mar-2020).
$ws.Range("C25").Value = "Shapefile for the Springield, Missouri MSA with populations interpolated to grid in field sum_pop_grid, sum_blk_grid, sum_lat_grid, sum_wht_grid (overwritten with corrected calculation on 10-Mar-2020; reference Homework01_workflow03.mxd)."

# The longer text now wraps onto two lines, so the row grows taller.
$ws.Range("A25:C25").RowHeight = 30

